$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the "Expected Behaviour" (validate definitions) column H4 ---
# Old validate4/validate6 used validate_Screenshot=...; new ones check the
# tabbar icon visibility instead (validate_isIconDisplayed).
$h4 = "validate1`n{`nvalidate_PageTitle=Compliance JS specs`n};`nvalidate2`n{`nvalidate_PageTitle=Native Tabbar JS Test`n};`nvalidate3`n{`nvalidate_Text_Exists=VT200-0576`n};`nvalidate4`n{`nvalidate_isIconDisplayed=tabbar_xpath,true`n};`nvalidate5`n{`nvalidate_Text_Exists=VT200-0578`n};`nvalidate6`n{`nvalidate_isIconDisplayed=tabbar_xpath,false`n};"
$ws.Range("H4").Value = $h4

# --- Update the "Description" (Steps) column G4 ---
# Old script had two TakeScreenshot() calls (VT200_0578_before / VT200_0578);
# new script replaces them with a single mid-script `validate4;` call and
# drops the final TakeScreenshot before validate6.
$g4 = "wait(5);`nvalidate1;`nlink_Click(tabbar_test_link);`nvalidate2;`nSelectTestToRun(VT200_0576_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(5);`nvalidate4;`nSwitchApp(NATIVE_APP);`nwait(2);`nClickNativeIcon(VT200_0576_mainpage_xpath);`nwait(2);`nSwitchApp(WEBVIEW);`nSelectTestToRun(VT200_0578_string);`nClickRunTest(runtest_top_xpath);`nvalidate5;`nClickRunTest(runtest_bottom_xpath);`nwait(5);`nvalidate6;"
$ws.Range("G4").Value = $g4

# --- Sheet view: scroll back to top and move the active selection to J1 ---
$ws.Range("J1").Select()
